
# Docx writer: Use different style for block quotes in notes.
# Add a new paragraph style "Footnote Block Text", based on "Footnote Text",
# mirroring the existing "Block Text" style's paragraph formatting
# (spacing + indents), so footnote block quotes can get their own font size.

$d = $word.ActiveDocument

$style = $d.Styles.Add("Footnote Block Text", 1)   # wdStyleTypeParagraph = 1

$style.BaseStyle = "Footnote Text"
$style.NextParagraphStyle = "Footnote Text"
$style.Priority = 9
$style.UnhideWhenUsed = $true
$style.QuickStyle = $true

$style.ParagraphFormat.SpaceBefore = 5
$style.ParagraphFormat.SpaceAfter = 5
$style.ParagraphFormat.FirstLineIndent = 0
$style.ParagraphFormat.LeftIndent = 24
$style.ParagraphFormat.RightIndent = 24

Write-Output ("Added style: " + $style.NameLocal)
